$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 0.655420508356823
$ws.Range("J3").Value = 0.5387434162345179
$ws.Range("K3").Value = 0.580364417305138
$ws.Range("L3").Value = 2.860943871911835
